# Applies the "2022 / Work in progress" update to Formats R3A 22.xlsx
# Updates the positional offsets (columns C "position" and D "fin") and
# lengths (column B "longueur") on sheet "Feuil1" for rows 14-29, to
# reflect newly inserted/lengthened fields, then moves the active
# selection to F29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 14
$ws.Range("B14").Value = 4
$ws.Range("D14").Value = 100

# Row 15
$ws.Range("C15").Value = 101
$ws.Range("D15").Value = 105

# Row 16
$ws.Range("C16").Value = 106
$ws.Range("D16").Value = 106

# Row 17
$ws.Range("C17").Value = 107
$ws.Range("D17").Value = 114

# Row 18
$ws.Range("C18").Value = 115
$ws.Range("D18").Value = 122

# Row 19
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = 123
$ws.Range("D19").Value = 124

# Row 20
$ws.Range("C20").Value = 125
$ws.Range("D20").Value = 127

# Row 21
$ws.Range("C21").Value = 128
$ws.Range("D21").Value = 128

# Row 22
$ws.Range("C22").Value = 129
$ws.Range("D22").Value = 129

# Row 23
$ws.Range("C23").Value = 130
$ws.Range("D23").Value = 130

# Row 24
$ws.Range("C24").Value = 131
$ws.Range("D24").Value = 131

# Row 25
$ws.Range("C25").Value = 132
$ws.Range("D25").Value = 132

# Row 26
$ws.Range("C26").Value = 133
$ws.Range("D26").Value = 133

# Row 27
$ws.Range("C27").Value = 134
$ws.Range("D27").Value = 139

# Row 28
$ws.Range("C28").Value = 140
$ws.Range("D28").Value = 141

# Row 29
$ws.Range("C29").Value = 142

# Move the active cell/selection to F29, matching the saved view state.
$ws.Activate()
$ws.Range("F29").Select() | Out-Null
